$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team names (column A) for the final (boys) point table ---
$ws.Range("A2").Value = "Mora A"
$ws.Range("A3").Value = "Sabra"
$ws.Range("A4").Value = "Pera"
$ws.Range("A5").Value = "Wayamba"
$ws.Range("A6").Value = "Rajarata"
$ws.Range("A7").Value = "Mora B"
$ws.Range("A8").Value = "Ruhuna"
$ws.Range("A9").Value = "Kelani"
$ws.Range("A10").Value = "Japura"
$ws.Range("A11").Value = "Colombo"

# --- Reset Played/Won/Drawn/Lost/Goals Scored/Goals Conceded to 0 for every team ---
$ws.Range("B2:G11").Value = 0

# --- Clear the penalty tally (Green/Yellow/Red) columns L:N entirely ---
$ws.Range("L2:N11").Clear()

# --- Column A width ---
$ws.Range("A:A").ColumnWidth = 13.1666667

# --- Selection as left by the author ---
$ws.Range("I15").Select()
